$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.6113651253405055
$ws.Range("J2").Value = 0.6113651253405055
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 4.712794644412445
$ws.Range("R2").Value = 42.415151799712
$ws.Range("S2").Value = 0.182387645985773
$ws.Range("T2").Value = 0.182387645985773

$ws.Range("I3").Value = 0.6113651253405055
$ws.Range("J3").Value = 0.6113651253405055
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("S3").Value = 0.1973402324458714
$ws.Range("T3").Value = 0.1973402324458714

$ws.Range("I4").Value = 0.6113651253405055
$ws.Range("J4").Value = 0.6113651253405055
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 2.198181498351889
$ws.Range("R4").Value = 19.783633485167
$ws.Range("S4").Value = 0.08507078690755557
$ws.Range("T4").Value = 0.08507078690755558

$ws.Range("I5").Value = 0.6113651253405055
$ws.Range("J5").Value = 0.6113651253405055
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 1.726171325088111
$ws.Range("R5").Value = 15.535541925793
$ws.Range("S5").Value = 0.06680374349097358
$ws.Range("T5").Value = 0.06680374349097358

$ws.Range("I6").Value = 0.6113651253405055
$ws.Range("J6").Value = 0.6113651253405055
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 2.061023931538667
$ws.Range("R6").Value = 18.549215383848
$ws.Range("S6").Value = 0.0797627165103319
$ws.Range("T6").Value = 0.07976271651033191

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01878466666666667
$ws.Range("H7").Value = 0.056354
$ws.Range("I7").Value = 0.3886348746594945
$ws.Range("J7").Value = 0.3886348746594945
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 2.995846966094223
$ws.Range("R7").Value = 26.962622694848
$ws.Range("S7").Value = 0.1159408625044529
$ws.Range("T7").Value = 0.1159408625044529

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.01878466666666667
$ws.Range("H8").Value = 0.056354
$ws.Range("I8").Value = 0.3886348746594945
$ws.Range("J8").Value = 0.3886348746594945
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("Q8").Value = 3.241453846646
$ws.Range("R8").Value = 29.173084619814
$ws.Range("S8").Value = 0.1254459787171565
$ws.Range("T8").Value = 0.1254459787171565

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.01878466666666667
$ws.Range("H9").Value = 0.056354
$ws.Range("I9").Value = 0.3886348746594945
$ws.Range("J9").Value = 0.3886348746594945
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 1.397348255046444
$ws.Range("R9").Value = 12.576134295418
$ws.Range("S9").Value = 0.05407811672049257
$ws.Range("T9").Value = 0.05407811672049258

$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.01878466666666667
$ws.Range("H10").Value = 0.056354
$ws.Range("I10").Value = 0.3886348746594945
$ws.Range("J10").Value = 0.3886348746594945
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 1.097299058713556
$ws.Range("R10").Value = 9.875691528422
$ws.Range("S10").Value = 0.04246605408501117
$ws.Range("T10").Value = 0.04246605408501117

$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.01878466666666667
$ws.Range("H11").Value = 0.056354
$ws.Range("I11").Value = 0.3886348746594945
$ws.Range("J11").Value = 0.3886348746594945
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 1.310159418821333
$ws.Range("R11").Value = 11.791434769392
$ws.Range("S11").Value = 0.05070386263238141
$ws.Range("T11").Value = 0.05070386263238142

